$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 0
$ws.Range("D10").Select() | Out-Null
